# Update NATMI LR-pair edge-weight metrics (Grn-Sort1) to the recomputed TPM values.
# Columns: G/H = ligand avg/total expression, I/J = ligand specificity (avg/total),
#          M/N = receptor avg/total expression, O/P = receptor specificity (avg/total),
#          Q/R = edge avg/total expression weight, S/T = edge specificity (avg/total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 32.00264033333334
$ws.Range("H2").Value = 96.00792100000001
$ws.Range("I2").Value = 0.02419001798940439
$ws.Range("J2").Value = 0.02433952891158457
$ws.Range("M2").Value = 0.4260053333333333
$ws.Range("N2").Value = 1.278016
$ws.Range("O2").Value = 0.02405532912416773
$ws.Range("P2").Value = 0.02531756756689831
$ws.Range("Q2").Value = 13.63329546274845
$ws.Range("R2").Value = 122.699659164736
$ws.Range("S2").Value = 0.0005818988442546605
$ws.Range("T2").Value = 0.0006162176677655172

# Row 3
$ws.Range("G3").Value = 32.00264033333334
$ws.Range("H3").Value = 96.00792100000001
$ws.Range("I3").Value = 0.02419001798940439
$ws.Range("J3").Value = 0.02433952891158457
$ws.Range("O3").Value = 0.05879323641880037
$ws.Range("P3").Value = 0.06187825274916518
$ws.Range("Q3").Value = 33.32091442903756
$ws.Range("R3").Value = 299.8882298613381
$ws.Range("S3").Value = 0.001422209446626086
$ws.Range("T3").Value = 0.001506087521786643

# Row 4
$ws.Range("G4").Value = 32.00264033333334
$ws.Range("H4").Value = 96.00792100000001
$ws.Range("I4").Value = 0.02419001798940439
$ws.Range("J4").Value = 0.02433952891158457
$ws.Range("M4").Value = 5.850740666666667
$ws.Range("N4").Value = 17.552222
$ws.Range("O4").Value = 0.3303749538898241
$ws.Range("P4").Value = 0.3477104875323931
$ws.Range("Q4").Value = 187.2391492389403
$ws.Range("R4").Value = 1685.152343150462
$ws.Range("S4").Value = 0.007991776077843489
$ws.Range("T4").Value = 0.008463109464155849

# Row 5
$ws.Range("G5").Value = 32.00264033333334
$ws.Range("H5").Value = 96.00792100000001
$ws.Range("I5").Value = 0.02419001798940439
$ws.Range("J5").Value = 0.02433952891158457
$ws.Range("M5").Value = 2.648771
$ws.Range("N5").Value = 5.297542
$ws.Range("O5").Value = 0.1495686865725097
$ws.Range("P5").Value = 0.1049445996947469
$ws.Range("Q5").Value = 84.76766563836368
$ws.Range("R5").Value = 508.6059938301821
$ws.Range("S5").Value = 0.003618069218840597
$ws.Range("T5").Value = 0.00255430211838496

# Row 6
$ws.Range("G6").Value = 32.00264033333334
$ws.Range("H6").Value = 96.00792100000001
$ws.Range("I6").Value = 0.02419001798940439
$ws.Range("J6").Value = 0.02433952891158457
$ws.Range("M6").Value = 7.742685666666667
$ws.Range("N6").Value = 23.228057
$ws.Range("O6").Value = 0.4372077939946981
$ws.Range("P6").Value = 0.4601490924567965
$ws.Range("Q6").Value = 247.7863846043886
$ws.Range("R6").Value = 2230.077461439497
$ws.Range("S6").Value = 0.01057606440183955
$ws.Range("T6").Value = 0.0111998121394916

# Row 7
$ws.Range("I7").Value = 0.07580145430919519
$ws.Range("J7").Value = 0.07626995934880827
$ws.Range("M7").Value = 0.4260053333333333
$ws.Range("N7").Value = 1.278016
$ws.Range("O7").Value = 0.02405532912416773
$ws.Range("P7").Value = 0.02531756756689831
$ws.Range("Q7").Value = 42.72107708046934
$ws.Range("R7").Value = 384.489693724224
$ws.Range("S7").Value = 0.001823428931498252
$ws.Range("T7").Value = 0.001930969849138041

# Row 8
$ws.Range("I8").Value = 0.07580145430919519
$ws.Range("J8").Value = 0.07626995934880827
$ws.Range("O8").Value = 0.05879323641880037
$ws.Range("P8").Value = 0.06187825274916518
$ws.Range("S8").Value = 0.004456612824089407
$ws.Range("T8").Value = 0.004719451821754112

# Row 9
$ws.Range("I9").Value = 0.07580145430919519
$ws.Range("J9").Value = 0.07626995934880827
$ws.Range("M9").Value = 5.850740666666667
$ws.Range("N9").Value = 17.552222
$ws.Range("O9").Value = 0.3303749538898241
$ws.Range("P9").Value = 0.3477104875323931
$ws.Range("Q9").Value = 586.7296097979287
$ws.Range("R9").Value = 5280.566488181358
$ws.Range("S9").Value = 0.02504290197218197
$ws.Range("T9").Value = 0.02651986474924993

# Row 10
$ws.Range("I10").Value = 0.07580145430919519
$ws.Range("J10").Value = 0.07626995934880827
$ws.Range("M10").Value = 2.648771
$ws.Range("N10").Value = 5.297542
$ws.Range("O10").Value = 0.1495686865725097
$ws.Range("P10").Value = 0.1049445996947469
$ws.Range("Q10").Value = 265.626604188473
$ws.Range("R10").Value = 1593.759625130838
$ws.Range("S10").Value = 0.01133752396131243
$ws.Range("T10").Value = 0.0080041203525953

# Row 11
$ws.Range("I11").Value = 0.07580145430919519
$ws.Range("J11").Value = 0.07626995934880827
$ws.Range("M11").Value = 7.742685666666667
$ws.Range("N11").Value = 23.228057
$ws.Range("O11").Value = 0.4372077939946981
$ws.Range("P11").Value = 0.4601490924567965
$ws.Range("Q11").Value = 776.4594602309636
$ws.Range("R11").Value = 6988.135142078672
$ws.Range("S11").Value = 0.03314098662011312
$ws.Range("T11").Value = 0.03509555257607089

# Row 12
$ws.Range("G12").Value = 473.968811
$ws.Range("H12").Value = 1421.906433
$ws.Range("I12").Value = 0.3582615042098434
$ws.Range("J12").Value = 0.360475806319893
$ws.Range("M12").Value = 0.4260053333333333
$ws.Range("N12").Value = 1.278016
$ws.Range("O12").Value = 0.02405532912416773
$ws.Range("P12").Value = 0.02531756756689831
$ws.Range("Q12").Value = 201.9132413196587
$ws.Range("R12").Value = 1817.219171876928
$ws.Range("S12").Value = 0.008618098396287185
$ws.Range("T12").Value = 0.00912637058273604

# Row 13
$ws.Range("G13").Value = 473.968811
$ws.Range("H13").Value = 1421.906433
$ws.Range("I13").Value = 0.3582615042098434
$ws.Range("J13").Value = 0.360475806319893
$ws.Range("O13").Value = 0.05879323641880037
$ws.Range("P13").Value = 0.06187825274916518
$ws.Range("Q13").Value = 493.4928502419194
$ws.Range("R13").Value = 4441.435652177274
$ws.Range("S13").Value = 0.02106335331676437
$ws.Range("T13").Value = 0.02230561305342146

# Row 14
$ws.Range("G14").Value = 473.968811
$ws.Range("H14").Value = 1421.906433
$ws.Range("I14").Value = 0.3582615042098434
$ws.Range("J14").Value = 0.360475806319893
$ws.Range("M14").Value = 5.850740666666667
$ws.Range("N14").Value = 17.552222
$ws.Range("O14").Value = 0.3303749538898241
$ws.Range("P14").Value = 0.3477104875323931
$ws.Range("Q14").Value = 2773.068597249347
$ws.Range("R14").Value = 24957.61737524413
$ws.Range("S14").Value = 0.1183606279338261
$ws.Range("T14").Value = 0.1253412183591225

# Row 15
$ws.Range("G15").Value = 473.968811
$ws.Range("H15").Value = 1421.906433
$ws.Range("I15").Value = 0.3582615042098434
$ws.Range("J15").Value = 0.360475806319893
$ws.Range("M15").Value = 2.648771
$ws.Range("N15").Value = 5.297542
$ws.Range("O15").Value = 0.1495686865725097
$ws.Range("P15").Value = 0.1049445996947469
$ws.Range("Q15").Value = 1255.434841481281
$ws.Range("R15").Value = 7532.609048887686
$ws.Range("S15").Value = 0.05358470263415795
$ws.Range("T15").Value = 0.03782998919388227

# Row 16
$ws.Range("G16").Value = 473.968811
$ws.Range("H16").Value = 1421.906433
$ws.Range("I16").Value = 0.3582615042098434
$ws.Range("J16").Value = 0.360475806319893
$ws.Range("M16").Value = 7.742685666666667
$ws.Range("N16").Value = 23.228057
$ws.Range("O16").Value = 0.4372077939946981
$ws.Range("P16").Value = 0.4601490924567965
$ws.Range("Q16").Value = 3669.791519376743
$ws.Range("R16").Value = 33028.12367439068
$ws.Range("S16").Value = 0.1566347219288079
$ws.Range("T16").Value = 0.1658726151307307

# Row 17
$ws.Range("G17").Value = 24.3798835
$ws.Range("H17").Value = 48.759767
$ws.Range("I17").Value = 0.01842816137361988
$ws.Range("J17").Value = 0.01236137337687614
$ws.Range("M17").Value = 0.4260053333333333
$ws.Range("N17").Value = 1.278016
$ws.Range("O17").Value = 0.02405532912416773
$ws.Range("P17").Value = 0.02531756756689831
$ws.Range("Q17").Value = 10.38596039704533
$ws.Range("R17").Value = 62.31576238227201
$ws.Range("S17").Value = 0.0004432954869957011
$ws.Range("T17").Value = 0.0003129599056887195

# Row 18
$ws.Range("G18").Value = 24.3798835
$ws.Range("H18").Value = 48.759767
$ws.Range("I18").Value = 0.01842816137361988
$ws.Range("J18").Value = 0.01236137337687614
$ws.Range("O18").Value = 0.05879323641880037
$ws.Range("P18").Value = 0.06187825274916518
$ws.Range("Q18").Value = 25.38415591438767
$ws.Range("R18").Value = 152.304935486326
$ws.Range("S18").Value = 0.001083451248403039
$ws.Range("T18").Value = 0.0007649001861411431

# Row 19
$ws.Range("G19").Value = 24.3798835
$ws.Range("H19").Value = 48.759767
$ws.Range("I19").Value = 0.01842816137361988
$ws.Range("J19").Value = 0.01236137337687614
$ws.Range("M19").Value = 5.850740666666667
$ws.Range("N19").Value = 17.552222
$ws.Range("O19").Value = 0.3303749538898241
$ws.Range("P19").Value = 0.3477104875323931
$ws.Range("Q19").Value = 142.6403758420457
$ws.Range("R19").Value = 855.8422550522741
$ws.Range("S19").Value = 0.006088202964083907
$ws.Range("T19").Value = 0.004298179163443546

# Row 20
$ws.Range("G20").Value = 24.3798835
$ws.Range("H20").Value = 48.759767
$ws.Range("I20").Value = 0.01842816137361988
$ws.Range("J20").Value = 0.01236137337687614
$ws.Range("M20").Value = 2.648771
$ws.Range("N20").Value = 5.297542
$ws.Range("O20").Value = 0.1495686865725097
$ws.Range("P20").Value = 0.1049445996947469
$ws.Range("Q20").Value = 64.5767283981785
$ws.Range("R20").Value = 258.306913592714
$ws.Range("S20").Value = 0.002756275892598583
$ws.Range("T20").Value = 0.001297259380713567

# Row 21
$ws.Range("G21").Value = 24.3798835
$ws.Range("H21").Value = 48.759767
$ws.Range("I21").Value = 0.01842816137361988
$ws.Range("J21").Value = 0.01236137337687614
$ws.Range("M21").Value = 7.742685666666667
$ws.Range("N21").Value = 23.228057
$ws.Range("O21").Value = 0.4372077939946981
$ws.Range("P21").Value = 0.4601490924567965
$ws.Range("Q21").Value = 188.7657745304532
$ws.Range("R21").Value = 1132.594647182719
$ws.Range("S21").Value = 0.008056935781538654
$ws.Range("T21").Value = 0.005688074740889159

# Row 22
$ws.Range("G22").Value = 692.3345543333334
$ws.Range("H22").Value = 2077.003663
$ws.Range("I22").Value = 0.5233188621179371
$ws.Range("J22").Value = 0.5265533320428379
$ws.Range("M22").Value = 0.4260053333333333
$ws.Range("N22").Value = 1.278016
$ws.Range("O22").Value = 0.02405532912416773
$ws.Range("P22").Value = 0.02531756756689831
$ws.Range("Q22").Value = 294.9382125969565
$ws.Range("R22").Value = 2654.443913372608
$ws.Range("S22").Value = 0.01258860746513193
$ws.Range("T22").Value = 0.01333104956156999

# Row 23
$ws.Range("G23").Value = 692.3345543333334
$ws.Range("H23").Value = 2077.003663
$ws.Range("I23").Value = 0.5233188621179371
$ws.Range("J23").Value = 0.5265533320428379
$ws.Range("O23").Value = 0.05879323641880037
$ws.Range("P23").Value = 0.06187825274916518
$ws.Range("Q23").Value = 720.8536608518016
$ws.Range("R23").Value = 6487.682947666214
$ws.Range("S23").Value = 0.03076760958291747
$ws.Range("T23").Value = 0.03258220016606182

# Row 24
$ws.Range("G24").Value = 692.3345543333334
$ws.Range("H24").Value = 2077.003663
$ws.Range("I24").Value = 0.5233188621179371
$ws.Range("J24").Value = 0.5265533320428379
$ws.Range("M24").Value = 5.850740666666667
$ws.Range("N24").Value = 17.552222
$ws.Range("O24").Value = 0.3303749538898241
$ws.Range("P24").Value = 0.3477104875323931
$ws.Range("Q24").Value = 4050.669931976577
$ws.Range("R24").Value = 36456.02938778919
$ws.Range("S24").Value = 0.1728914449418887
$ws.Range("T24").Value = 0.1830881157964213

# Row 25
$ws.Range("G25").Value = 692.3345543333334
$ws.Range("H25").Value = 2077.003663
$ws.Range("I25").Value = 0.5233188621179371
$ws.Range("J25").Value = 0.5265533320428379
$ws.Range("M25").Value = 2.648771
$ws.Range("N25").Value = 5.297542
$ws.Range("O25").Value = 0.1495686865725097
$ws.Range("P25").Value = 0.1049445996947469
$ws.Range("Q25").Value = 1833.835689816058
$ws.Range("R25").Value = 11003.01413889635
$ws.Range("S25").Value = 0.07827211486560018
$ws.Range("T25").Value = 0.05525892864917075

# Row 26 (source diff hunk for this row was truncated/corrupted - missing the
# Q26/R26/S26 lines and showing a garbled old/new pair for P26. Values below
# reconstructed from the same Q=G*M, R=H*N, S=Q/sum(Q), T=R/sum(R) relationship
# that holds for every other row, and P26 follows the identical M/N/O/P group
# pattern shared with rows 6/11/16/21; the recovered T26 matches the diff.)
$ws.Range("G26").Value = 692.3345543333334
$ws.Range("H26").Value = 2077.003663
$ws.Range("I26").Value = 0.5233188621179371
$ws.Range("J26").Value = 0.5265533320428379
$ws.Range("M26").Value = 7.742685666666667
$ws.Range("N26").Value = 23.228057
$ws.Range("O26").Value = 0.4372077939946981
$ws.Range("P26").Value = 0.4601490924567965
$ws.Range("Q26").Value = 5360.528830374755
$ws.Range("R26").Value = 48244.75947337279
$ws.Range("S26").Value = 0.2287990852623988
$ws.Range("T26").Value = 0.2422930378696141
